# Update build timestamp strings throughout the workbook.
# Old build timestamp: "January 30 2026 16.19.47 EST"
# New build timestamp: "February 02 2026 12.49.33 EST"

$wb = $excel.ActiveWorkbook

$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

$aboutSheet = $wb.Worksheets.Item("About")
$dataSheet  = $wb.Worksheets.Item("Boundaries and methane sources")

# A2: "Version: mines - January 30 (built on ... )"
$a2 = $aboutSheet.Range("A2")
$a2.Value = $a2.Value().Replace($oldStamp, $newStamp)

# A6: Recommended citation text
$a6 = $aboutSheet.Range("A6")
$a6.Value = $a6.Value().Replace($oldStamp, $newStamp)

# S2:S9 on the data sheet: build_version column
for ($row = 2; $row -le 9; $row++) {
    $cell = $dataSheet.Cells.Item($row, 19)  # column S = 19
    $cell.Value = $cell.Value().Replace($oldStamp, $newStamp)
}
